$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D2..D48: numeric-looking price strings. A leading apostrophe forces
# Excel to keep them as literal text instead of coercing to a Double; the
# explicit Style reset afterward clears the "quote prefix" flag that the
# apostrophe entry sets, so no stray style/number-format is left behind.

$ws.Range("D2").Value = "'242.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.378"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05927"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.397"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.443"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8062"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9062"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1419"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07424"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03229"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09309"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.952"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001577"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04772"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006132"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006136"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "UpBots"
$ws.Range("C20").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D20").Value = "'0.007493"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19UpBotsUBXTBestin24h"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.004391"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.0009807"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.00007805"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.611"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.150"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("D26").Value = "'0.3249"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.1323"
$ws.Range("D27").Style = "Normal"
$ws.Range("D41").Value = "'0.006211"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1062"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002621"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007264"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005201"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.0005811"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.9619"
$ws.Range("D48").Style = "Normal"
